$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(91, 1).Value = 'Clase 3'
$ws.Cells.Item(92, 1).Value = 'necesito hacerlo'
$ws.Cells.Item(92, 2).Value = 'j''ai besoin de le faire'
$ws.Cells.Item(93, 1).Value = 'no necesito hacerlo'
$ws.Cells.Item(93, 2).Value = 'je n''ai pas besoin de le faire'
$ws.Cells.Item(94, 1).Value = 'necesito estar aquí'
$ws.Cells.Item(94, 2).Value = 'j''ai besoin d''être ici'
$ws.Cells.Item(95, 1).Value = 'no necesito estar aquí'
$ws.Cells.Item(95, 2).Value = 'je n''ai pas besoin d''être ici'
$ws.Cells.Item(96, 1).Value = 'no necesito estar aquí hoy'
$ws.Cells.Item(96, 2).Value = 'je n''ai pas besoin d''être ici aujourd''hui'
$ws.Cells.Item(97, 1).Value = 'necesito hacerlo ahora'
$ws.Cells.Item(97, 2).Value = 'j''ai besoin d''être maintenant'
$ws.Cells.Item(98, 1).Value = 'necesito estar aquí ahora'
$ws.Cells.Item(98, 2).Value = 'j''ai besoin d''être ici maintenant'
$ws.Cells.Item(99, 1).Value = 'porque'
$ws.Cells.Item(99, 2).Value = 'parce que'
$ws.Cells.Item(100, 1).Value = 'porque es importante'
$ws.Cells.Item(100, 2).Value = 'parce que c''est très important'
$ws.Cells.Item(101, 1).Value = 'necesito hacerlo hoy porque es importante'
$ws.Cells.Item(101, 2).Value = 'j''ai besoin de le faire aujourd''hui parce que c''est important'
$ws.Cells.Item(102, 1).Value = 'necesito estar aquí hoy porque es importante'
$ws.Cells.Item(102, 2).Value = 'j''ai besoin d''être ici aujourd''hui parce que c''est important'
$ws.Cells.Item(103, 1).Value = 'es posible'
$ws.Cells.Item(103, 2).Value = 'c''est possible'
$ws.Cells.Item(104, 1).Value = 'no es posible'
$ws.Cells.Item(104, 2).Value = 'ce n''est pas possible'
$ws.Cells.Item(105, 1).Value = 'ahora, no es posible'
$ws.Cells.Item(105, 2).Value = 'maintenant, ce n''est pas possible'
$ws.Cells.Item(106, 1).Value = 'saber'
$ws.Cells.Item(106, 2).Value = 'savoir'
$ws.Cells.Item(107, 1).Value = 'necesito saber'
$ws.Cells.Item(107, 2).Value = 'j''ai besoin de savoir'
$ws.Cells.Item(108, 1).Value = 'si (condicional)'
$ws.Cells.Item(108, 2).Value = 'si'
$ws.Cells.Item(109, 1).Value = 'sí (afirmacion)'
$ws.Cells.Item(109, 2).Value = 'oui'
$ws.Cells.Item(110, 1).Value = 'si es posible'
$ws.Cells.Item(110, 2).Value = 'si c''est possible'
$ws.Cells.Item(111, 1).Value = 'si es posible para mi'
$ws.Cells.Item(111, 2).Value = 'si c''est possible pour moi'
$ws.Cells.Item(112, 1).Value = 'necesito saber si es posible'
$ws.Cells.Item(112, 2).Value = 'j''ai besoin de savoir si c''est possible'
$ws.Cells.Item(113, 1).Value = 'necesito saber si es posible hoy'
$ws.Cells.Item(113, 2).Value = 'j''ai besoin de savoir si c''est possible aujourd''hui'
$ws.Cells.Item(114, 1).Value = 'necesito saber si es posible hacerlo hoy'
$ws.Cells.Item(114, 2).Value = 'j''ai besoin de savoir si c''est possible de le faire aujourd''hui'
$ws.Cells.Item(115, 1).Value = 'necesito saber si es posible hacerlo ahora porque es muy importante'
$ws.Cells.Item(115, 2).Value = 'j''ai besoin de savoir si c''est possible de le faire maintenant parce que c''est très important'
$ws.Cells.Item(116, 1).Value = 'necesito saber si es algo importante'
$ws.Cells.Item(116, 2).Value = 'j''ai besoin de savoir si c''est quelque chose d''important'
$ws.Cells.Item(117, 1).Value = 'yo quiero'
$ws.Cells.Item(117, 2).Value = 'je veux'
$ws.Cells.Item(118, 1).Value = 'quiero hacerlo'
$ws.Cells.Item(118, 2).Value = 'je veux le faire'
$ws.Cells.Item(119, 1).Value = 'quiero hacerlo hoy'
$ws.Cells.Item(119, 2).Value = 'je veux le faire aujourd''hui'
$ws.Cells.Item(120, 1).Value = 'no quiero'
$ws.Cells.Item(120, 2).Value = 'je ne veux pas'
$ws.Cells.Item(121, 1).Value = 'no quiero hacerlo así'
$ws.Cells.Item(121, 2).Value = 'je ne veux pas le faire comme ça'
$ws.Cells.Item(122, 1).Value = 'quiero saber si es fácil para mí hacerlo hoy'
$ws.Cells.Item(122, 2).Value = 'je veux savoir si c''est facile pour moi de le faire aujourd''hui'

$ws.Columns.Item(1).ColumnWidth = 30.71
$ws.Columns.Item(2).ColumnWidth = 19.29

$ws.Range("A97").Select()
$excel.ActiveWindow.ScrollRow = 97
$ws.Range("B123").Select()
